$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data in rows 3 and 4 for the columns that differ between them:
# D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), Q (Unidad de comercializacion),
# S (Precio $/Kg), T (Kg / unidad)

$ws.Range("D3").Value = 44875
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 16000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("S3").Value = 1600
$ws.Range("T3").Value = 10

$ws.Range("D4").Value = 44855
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "$/bandeja 5 kilos"
$ws.Range("S4").Value = 3000
$ws.Range("T4").Value = 5
